# Table 16: Subjective well-being — refresh to the 2012-21 measurement window.
# (OECD "Add files via upload" commit: updates header captions, refreshes
# the underlying indicator figures for a batch of countries, and marks
# Djibouti's row as unavailable ("..") for this edition.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab16")

# --- Header captions (row 2, columns C:H): "2011-20" -> "2012-21" ---
$ws.Range("C2").Value2 = "Cantril life ladder, most recent measure 2012-21"
$ws.Range("D2").Value2 = "Negative affect, most recent measure 2012-21"
$ws.Range("E2").Value2 = "Positive affect, most recent measure 2012-21"
$ws.Range("F2").Value2 = "Freedom to make life choices, most recent measure 2012-21"
$ws.Range("G2").Value2 = "Generosity, most recent measure 2012-21"
$ws.Range("H2").Value2 = "Social support, most recent measure 2012-21"

# Djibouti (row 25): all six indicators now unavailable ("..")
$ws.Range("C25:H25").Value2 = ".."

# Refreshed indicator values (Table 16, 2012-21 measurement window)
# Row 38
$ws.Range("C38").Value2 = 4.3069809350100403
$ws.Range("D38").Value2 = 0.31398930874738001
$ws.Range("E38").Value2 = 0.70614951307122997
$ws.Range("F38").Value2 = 0.68648891557347003
$ws.Range("G38").Value2 = 0.10217988242706
$ws.Range("H38").Value2 = 0.70123958858576996
# Row 62
$ws.Range("C62").Value2 = 4.5010495896035101
$ws.Range("D62").Value2 = 0.34835901127217
$ws.Range("E62").Value2 = 0.67362115738240003
$ws.Range("F62").Value2 = 0.69406166419070003
$ws.Range("G62").Value2 = 0.0077786147490999996
$ws.Range("H62").Value2 = 0.68738903327191003
# Row 63
$ws.Range("C63").Value2 = 5.8937654735845202
$ws.Range("D63").Value2 = 0.29445539704627
$ws.Range("F63").Value2 = 0.82816123142154996
$ws.Range("G63").Value2 = -0.0047972306924999996
# Row 66
$ws.Range("C66").Value2 = 5.47416517520562
$ws.Range("D66").Value2 = 0.31080036394058003
$ws.Range("E66").Value2 = 0.70677966129394998
$ws.Range("F66").Value2 = 0.78775943873020005
$ws.Range("G66").Value2 = -0.0011224706609000001
$ws.Range("H66").Value2 = 0.80146210678876995
# Row 67
$ws.Range("C67").Value2 = 4.4254035022523697
$ws.Range("D67").Value2 = 0.33511081337929
$ws.Range("E67").Value2 = 0.68409664101070999
$ws.Range("F67").Value2 = 0.69834480186303005
$ws.Range("G67").Value2 = 0.0083654248703500002
$ws.Range("H67").Value2 = 0.69900141656398995
# Row 68
$ws.Range("C68").Value2 = 4.6771178245544398
$ws.Range("D68").Value2 = 0.36442494595593
$ws.Range("E68").Value2 = 0.67578378319739996
$ws.Range("F68").Value2 = 0.70023567703637002
$ws.Range("G68").Value2 = 0.01269218447574
$ws.Range("H68").Value2 = 0.65311868895184
# Row 72
$ws.Range("C72").Value2 = 4.2266583045323696
$ws.Range("D72").Value2 = 0.33109540492296002
$ws.Range("E72").Value2 = 0.68660712242125999
$ws.Range("F72").Value2 = 0.65360676248867999
$ws.Range("G72").Value2 = 0.13304734043777
$ws.Range("H72").Value2 = 0.70986185471217
# Row 81
$ws.Range("C81").Value2 = 5.6594903048346996
$ws.Range("D81").Value2 = 0.28630470440667999
$ws.Range("F81").Value2 = 0.81030523075776995
$ws.Range("G81").Value2 = 0.0034878488410900002
# Row 82
$ws.Range("C82").Value2 = 4.4524682118342502
$ws.Range("D82").Value2 = 0.34706117403813003
$ws.Range("E82").Value2 = 0.67761577704013998
$ws.Range("F82").Value2 = 0.71278704282564997
$ws.Range("G82").Value2 = 0.021118091822070001
$ws.Range("H82").Value2 = 0.67582122561258995
# Row 86
$ws.Range("C86").Value2 = 4.5986801385879499
$ws.Range("D86").Value2 = 0.32320380881429001
$ws.Range("E86").Value2 = 0.67545933127403002
$ws.Range("F86").Value2 = 0.68152459114789998
$ws.Range("G86").Value2 = -0.0116575255292
$ws.Range("H86").Value2 = 0.69564319849014
# Row 90
$ws.Range("C90").Value2 = 6.6837128003438302
$ws.Range("D90").Value2 = 0.25683281852139001
$ws.Range("F90").Value2 = 0.86471313767962998
$ws.Range("G90").Value2 = -0.025535291112500001
# Row 91
$ws.Range("C91").Value2 = 4.3115492853625099
$ws.Range("D91").Value2 = 0.36660103499889002
$ws.Range("E91").Value2 = 0.66614309056052001
$ws.Range("F91").Value2 = 0.68185059880387999
$ws.Range("G91").Value2 = 0.043668594797729998
$ws.Range("H91").Value2 = 0.66145415244431005
# Row 97
$ws.Range("C97").Value2 = 4.4438498959396799
$ws.Range("D97").Value2 = 0.36362232448477
$ws.Range("E97").Value2 = 0.67520290974414998
$ws.Range("F97").Value2 = 0.67108643687132996
$ws.Range("G97").Value2 = 0.042571862423490002
$ws.Range("H97").Value2 = 0.68556041609157004
